# Adds a new task row (row 11) to the Tarefas sheet:
#   A11 = "Criar verificação na leitura dos valores das arestas onde não pode ser negativa"
#   B11 = "Douglas"
# and moves the selection / top-left cell back to A1 / D11 as in the author's
# saved view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = "Criar verificação na leitura dos valores das arestas onde não pode ser negativa"
$ws.Range("B11").Value = "Douglas"

$ws.Range("A1").Select()
$ws.Range("D11").Select()
